# DDAf_2022_Tableau_annexe_Tab16.xlsx - text corrections on the "Tab16" sheet
# (country-name label cleanups in column B) plus the window-size bookkeeping
# change recorded in the workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab16")

# Window size captured in <bookViews><workbookView .../></bookViews> when the
# author last saved the file (best effort - harmless if the host ignores it).
$excel.ActiveWindow.Height = 12490

# Column B country / territory label corrections (shared strings).
$ws.Range("B4").Value  = "Botswana"                 # was "Botswana*"
$ws.Range("B11").Value = "Zambie"                   # was "Zambie*"
$ws.Range("B18").Value = "République du Congo*"     # was "Congo*"
$ws.Range("B19").Value = "RD Congo"                 # was "DR Congo*"
$ws.Range("B22").Value = "Sao Tomé-et-Principe"      # was "São Tomé and Príncipe"
$ws.Range("B34").Value = "Soudan du Sud"            # was "Soudan du Sud*"
$ws.Range("B36").Value = "Tanzanie"                 # was "UR of Tanzania"
$ws.Range("B54").Value = "Liberia"                  # was "Libéria"
$ws.Range("B57").Value = "Nigeria*"                 # was "Nigéria*"
$ws.Range("B60").Value = "Togo"                     # was "Togo*"
